# Apply updated "想去人数" (F column) counts to the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$changes = @(
    @{Row = 6;  Old = 2999; New = 3001},
    @{Row = 8;  Old = 2016; New = 2019},
    @{Row = 9;  Old = 320;  New = 321},
    @{Row = 12; Old = 954;  New = 955},
    @{Row = 14; Old = 421;  New = 424},
    @{Row = 19; Old = 7260; New = 7273},
    @{Row = 21; Old = 1985; New = 1994},
    @{Row = 22; Old = 190;  New = 191},
    @{Row = 27; Old = 75;   New = 74},
    @{Row = 28; Old = 1121; New = 1122},
    @{Row = 40; Old = 151;  New = 152},
    @{Row = 41; Old = 279;  New = 280}
)

$changes4 = @(
    @{Row = 9;  Old = 2999; New = 3001},
    @{Row = 11; Old = 2016; New = 2019},
    @{Row = 12; Old = 320;  New = 321},
    @{Row = 16; Old = 954;  New = 955},
    @{Row = 18; Old = 421;  New = 424},
    @{Row = 23; Old = 7260; New = 7273},
    @{Row = 25; Old = 1985; New = 1994},
    @{Row = 27; Old = 190;  New = 191},
    @{Row = 32; Old = 75;   New = 74},
    @{Row = 33; Old = 1121; New = 1122},
    @{Row = 44; Old = 151;  New = 152},
    @{Row = 45; Old = 279;  New = 280}
)

$ws1 = $wb.Worksheets.Item("展览")
foreach ($chg in $changes) {
    $cell = $ws1.Range("F" + $chg.Row)
    if ($cell.Value() -eq $chg.Old) {
        $cell.Value = $chg.New
    }
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($chg in $changes4) {
    $cell = $ws4.Range("F" + $chg.Row)
    if ($cell.Value() -eq $chg.Old) {
        $cell.Value = $chg.New
    }
}
